$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "First available time" value for row 3 (barge #2)
$ws.Range("I3").Value = 201801010615

# Move the active selection to K1 (as recorded in the saved view state)
$ws.Range("K1").Select()
